$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.593.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = "'3.687.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'665.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").Value = "'160.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = "'0.499"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = "'7.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.23%  '
$ws.Range("D11").Value = "'0.441"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D13").Value = "'32.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = "'3.691.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = "'69.561.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("E16").Value = '  +2.68%  '
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").Value = "'6.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").Value = "'469.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = "'9.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = "'0.646"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").Value = "'79.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.833.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = "'0.0000127"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'10.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = "'2.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").Value = '  -3.37%  '
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'26.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = "'0.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.10%  '
$ws.Range("D34").Value = "'6.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("D35").Value = "'3.678.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = "'8.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.06%  '
$ws.Range("D39").Value = "'179.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.59%  '
$ws.Range("D40").Value = "'2.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = "'0.934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("D44").Value = "'47.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").Value = "'2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.81%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = "'27.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = "'0.000270"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.11%  '
$ws.Range("D49").Value = "'1.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").Value = "'7.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").Value = "'0.264"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.68%  '
